$wb = $excel.ActiveWorkbook
Write-Host ($wb.Worksheets | ForEach-Object { $_.Name }) -join ","
